# Applies the data refresh from the commit "Add files via upload" to
# Data/日间手术总结报表.xlsx — updates cumulative flow counters across the
# workbook's twelve sheets (new month/week of data appended + totals bumped).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 日间手术年流量 (sheet1) — yearly total
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("日间手术年流量")
$ws.Range("B2").Value = 817

# ---------------------------------------------------------------------
# 2. 月流量 (sheet2) — monthly counts, month 10 revised + month 11 added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("月流量")
$ws.Cells.Item(11, 2).Value = 95
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 125

# ---------------------------------------------------------------------
# 3. 周流量 (sheet3) — weekly counts, week 44 revised + weeks 45-48 added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("周流量")
$ws.Cells.Item(38, 2).Value = 23
$ws.Cells.Item(39, 1).Value = 45
$ws.Cells.Item(39, 2).Value = 33
$ws.Cells.Item(40, 1).Value = 46
$ws.Cells.Item(40, 2).Value = 34
$ws.Cells.Item(41, 1).Value = 47
$ws.Cells.Item(41, 2).Value = 31
$ws.Cells.Item(42, 1).Value = 48
$ws.Cells.Item(42, 2).Value = 23

# ---------------------------------------------------------------------
# 4. 周内流量 (sheet4) — weekday-of-week counts, rows 2-8
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("周内流量")
$ws.Cells.Item(2, 2).Value = 163
$ws.Cells.Item(3, 2).Value = 157
$ws.Cells.Item(4, 2).Value = 157
$ws.Cells.Item(5, 2).Value = 118
$ws.Cells.Item(6, 2).Value = 139
$ws.Cells.Item(7, 2).Value = 41
$ws.Cells.Item(8, 2).Value = 42

# ---------------------------------------------------------------------
# 5. 核算年 (sheet5) — yearly total (accounting view)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("核算年")
$ws.Range("B2").Value = 817

# ---------------------------------------------------------------------
# 6. 核算月 (sheet6) — month 11 revised + month 12 added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("核算月")
$ws.Cells.Item(11, 2).Value = 128
$ws.Cells.Item(12, 1).Value = 12
$ws.Cells.Item(12, 2).Value = 18

# ---------------------------------------------------------------------
# 7. 核算年核算月流量 (sheet7) — month 11 revised + month 12 column added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("核算年核算月流量")
$ws.Cells.Item(1, 12).Value = 12
$ws.Cells.Item(2, 11).Value = 128
$ws.Cells.Item(2, 12).Value = 18

# ---------------------------------------------------------------------
# 8. 年周期月度流量 (sheet8) — month 10 revised + month 11 column added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("年周期月度流量")
$ws.Cells.Item(1, 12).Value = 11
$ws.Cells.Item(2, 11).Value = 95
$ws.Cells.Item(2, 12).Value = 125

# ---------------------------------------------------------------------
# 9. 年周期周度流量 (sheet9) — week 44 revised + weeks 45-48 columns added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("年周期周度流量")
$ws.Cells.Item(1, 39).Value = 45
$ws.Cells.Item(1, 40).Value = 46
$ws.Cells.Item(1, 41).Value = 47
$ws.Cells.Item(1, 42).Value = 48
$ws.Cells.Item(2, 38).Value = 23
$ws.Cells.Item(2, 39).Value = 33
$ws.Cells.Item(2, 40).Value = 34
$ws.Cells.Item(2, 41).Value = 31
$ws.Cells.Item(2, 42).Value = 23

# ---------------------------------------------------------------------
# 10. 年周期月中流量 (sheet10) — every day-of-month column revised
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("年周期月中流量")
$ws.Cells.Item(2, 2).Value = 12
$ws.Cells.Item(2, 3).Value = 14
$ws.Cells.Item(2, 4).Value = 17
$ws.Cells.Item(2, 5).Value = 25
$ws.Cells.Item(2, 6).Value = 20
$ws.Cells.Item(2, 7).Value = 26
$ws.Cells.Item(2, 8).Value = 29
$ws.Cells.Item(2, 9).Value = 29
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(2, 11).Value = 22
$ws.Cells.Item(2, 12).Value = 38
$ws.Cells.Item(2, 13).Value = 36
$ws.Cells.Item(2, 14).Value = 32
$ws.Cells.Item(2, 15).Value = 30
$ws.Cells.Item(2, 16).Value = 31
$ws.Cells.Item(2, 17).Value = 22
$ws.Cells.Item(2, 18).Value = 21
$ws.Cells.Item(2, 19).Value = 28
$ws.Cells.Item(2, 20).Value = 31
$ws.Cells.Item(2, 21).Value = 37
$ws.Cells.Item(2, 22).Value = 30
$ws.Cells.Item(2, 23).Value = 32
$ws.Cells.Item(2, 24).Value = 21
$ws.Cells.Item(2, 25).Value = 34
$ws.Cells.Item(2, 26).Value = 28
$ws.Cells.Item(2, 27).Value = 31
$ws.Cells.Item(2, 28).Value = 24
$ws.Cells.Item(2, 29).Value = 26
$ws.Cells.Item(2, 30).Value = 33
$ws.Cells.Item(2, 31).Value = 20
$ws.Cells.Item(2, 32).Value = 13

# ---------------------------------------------------------------------
# 11. 年周期周中流量 (sheet11) — weekday-of-week counts (single row)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("年周期周中流量")
$ws.Cells.Item(2, 2).Value = 163
$ws.Cells.Item(2, 3).Value = 157
$ws.Cells.Item(2, 4).Value = 157
$ws.Cells.Item(2, 5).Value = 118
$ws.Cells.Item(2, 6).Value = 139
$ws.Cells.Item(2, 7).Value = 41
$ws.Cells.Item(2, 8).Value = 42

# ---------------------------------------------------------------------
# 12. 月周期周中流量 (sheet12) — week-of-month 10 revised + week 11 added
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("月周期周中流量")
$ws.Cells.Item(11, 5).Value = 14
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 28
$ws.Cells.Item(12, 3).Value = 24
$ws.Cells.Item(12, 4).Value = 19
$ws.Cells.Item(12, 5).Value = 17
$ws.Cells.Item(12, 6).Value = 22
$ws.Cells.Item(12, 7).Value = 7
$ws.Cells.Item(12, 8).Value = 8
